# "Data updated 13/4 12:00"
# Corrects several historical D/G values and appends the 5 new daily rows
# (one per province) for 2020-04-13 (serial 43934).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Keep the existing autofilter range pinned to the original data ----
# (must run before the new rows are appended, otherwise Excel grows the
# filter to cover the whole used range automatically)
$ws.Range("A1:G136").AutoFilter() | Out-Null
$fdb = $ws.Names.Add("_xlnm._FilterDatabase", "=Sheet1!`$A`$1:`$G`$136")
$fdb.Visible = $false

# --- Corrections to previously-entered historical values ---------------
$corrections = @(
    @{ Row = 32;  Col = "D"; Value = 5 }
    @{ Row = 36;  Col = "D"; Value = 13 }
    @{ Row = 41;  Col = "D"; Value = 14 }
    @{ Row = 53;  Col = "D"; Value = 3 }
    @{ Row = 56;  Col = "D"; Value = 45 }
    @{ Row = 61;  Col = "D"; Value = 45 }
    @{ Row = 62;  Col = "D"; Value = 15 }
    @{ Row = 67;  Col = "D"; Value = 14 }
    @{ Row = 72;  Col = "D"; Value = 6 }
    @{ Row = 76;  Col = "D"; Value = 27 }
    @{ Row = 77;  Col = "G"; Value = 1 }
    @{ Row = 78;  Col = "D"; Value = 31 }
    @{ Row = 81;  Col = "D"; Value = 111 }
    @{ Row = 86;  Col = "D"; Value = 123 }
    @{ Row = 87;  Col = "D"; Value = 53 }
    @{ Row = 91;  Col = "D"; Value = 62 }
    @{ Row = 92;  Col = "D"; Value = 21 }
    @{ Row = 92;  Col = "G"; Value = 2 }
    @{ Row = 96;  Col = "D"; Value = 80 }
    @{ Row = 96;  Col = "G"; Value = 0 }
)

foreach ($fix in $corrections) {
    $ws.Range($fix.Col + $fix.Row).Value = $fix.Value
}

# --- Append the five new rows for 2020-04-13 (serial 43934) -------------
$newRows = @(
    @{ Id = 135; Date = 43934; Provincia = "LAS TUNAS";   D = 66; E = 0; F = 0; G = 0 }
    @{ Id = 136; Date = 43934; Provincia = "HOLGUÍN";     D = 22; E = 0; F = 1; G = 0 }
    @{ Id = 137; Date = 43934; Provincia = "GRANMA";      D = 17; E = 0; F = 0; G = 0 }
    @{ Id = 138; Date = 43934; Provincia = "SANTIAGO";    D = 39; E = 0; F = 0; G = 0 }
    @{ Id = 139; Date = 43934; Provincia = "GUANTÁNAMO";  D = 40; E = 0; F = 2; G = 0 }
)

$startRow = 137
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $item = $newRows[$i]

    $ws.Range("A$r").Value = $item.Id
    $ws.Range("B$r").Value = $item.Date
    $ws.Range("C$r").Value = $item.Provincia
    $ws.Range("D$r").Value = $item.D
    $ws.Range("E$r").Value = $item.E
    $ws.Range("F$r").Value = $item.F
    $ws.Range("G$r").Value = $item.G

    # Match the formatting used by the existing data rows (border/bold
    # id style in col A, date style in col B, plain in C:G).
    $ws.Range("A136").Copy() | Out-Null
    $ws.Range("A$r").PasteSpecial(-4122) | Out-Null
    $ws.Range("B136").Copy() | Out-Null
    $ws.Range("B$r").PasteSpecial(-4122) | Out-Null
    $ws.Range("C136:G136").Copy() | Out-Null
    $ws.Range("C$r:G$r").PasteSpecial(-4122) | Out-Null
}
$excel.CutCopyMode = $false

# --- Re-freeze header row / first column, keep the selection where the
# user last worked (bottom-right pane at the newly-added cell) ----------
$win = $excel.ActiveWindow
$win.FreezePanes = $false
$ws.Range("B2").Select() | Out-Null
$win.FreezePanes = $true
$ws.Range("D141").Select() | Out-Null

# --- Page setup / orientation, same as before, now explicit ------------
$ws.PageSetup.Orientation = 1
